{"js": "// UC3 \"BeregnPrioritet\" use case \u2014 dress up the casual description:\n//  1. Title \"BeregnPrioritet\" -> \"Beregn prioritet\"\n//  2. Extension step \"Risk Manager udskriver fejl besked\" -> \"Risk Manager viser fejl besked\"\n\n// --- 1) Fix the (single-word) title into two properly spaced words ---\nconst titleResults = context.document.body.search(\"BeregnPrioritet\", { matchCase: true, matchWholeWord: false });\ntitleResults.load(\"text\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\"Beregn prioritet\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2) Swap \"udskriver\" for \"viser\" in the error-message step, keeping the\n//        trailing \" og stopper udregning\" (held in a separate run) intact ---\nconst stepResults = context.document.body.search(\"Risk Manager udskriver fejl besked\", { matchCase: true, matchWholeWord: false });\nstepResults.load(\"text\");\nawait context.sync();\n\nif (stepResults.items.length > 0) {\n  stepResults.items[0].insertText(\"Risk Manager viser fejl besked\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# UC3 \"BeregnPrioritet\" use case \u2014 dress up the casual description:\n#  1. Title \"BeregnPrioritet\" -> \"Beregn prioritet\"\n#  2. Extension step \"Risk Manager udskriver fejl besked\" -> \"Risk Manager viser fejl besked\"\n\n$d = $word.ActiveDocument\n\n# --- 1) Fix the (single-word) title into two properly spaced words ---\n$titleFind = $d.Content.Find\n$titleFind.ClearFormatting()\n$titleFind.Replacement.ClearFormatting()\n$titleFind.Text = \"BeregnPrioritet\"\n$titleFind.Replacement.Text = \"Beregn prioritet\"\n$titleFind.Execute([ref]$titleFind.Text, [ref]$true, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$true, 1, [ref]$false, [ref]$titleFind.Replacement.Text, 2) | Out-Null\n\n# --- 2) Swap \"udskriver\" for \"viser\" in the error-message step, keeping the\n#        trailing \" og stopper udregning\" text intact ---\n$stepFind = $d.Content.Find\n$stepFind.ClearFormatting()\n$stepFind.Replacement.ClearFormatting()\n$stepFind.Text = \"Risk Manager udskriver fejl besked\"\n$stepFind.Replacement.Text = \"Risk Manager viser fejl besked\"\n$stepFind.Execute([ref]$stepFind.Text, [ref]$true, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$true, 1, [ref]$false, [ref]$stepFind.Replacement.Text, 2) | Out-Null\n"}
